$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.303.41'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '1.608.58'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''212.77'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = '''18.51'
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = '1.831.32'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '1.606.21'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '''4.03'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '''0.516'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '26.270.24'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '''62.20'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '''201.35'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = '''4.27'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").Value = '''143.48'
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '''6.58'
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("E30").Value = '  +4.87%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D36").Value = '1.163.49'
$ws.Range("E36").Value = '  +3.30%  '
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '''2.32'
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("D42").Value = '''5.38'
$ws.Range("E42").Value = '  +4.34%  '
$ws.Range("D43").Value = '''0.786'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '1.742.63'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '''92.11'
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''1.54'
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0105'
$ws.Range("E47").Value = '  +13.29%  '
$ws.Range("D48").Value = '''54.12'
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("E51").Value = '  -0.03%  '
